# Weekly update: insert the latest week's "Choclo" (Feria Lagunitas de
# Puerto Montt) record. This shifts every existing data row down by one
# (row 182 -> 183, ..., old row 241 -> new row 242) and fills the freed
# row 182 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 182..241 down to 183..242, opening up a blank row 182.
$ws.Rows.Item(182).Insert()

# Populate the newly opened row 182 with the new record.
$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = 44627
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = 100112024
$ws.Range("G182").Value = "Choclo"
$ws.Range("H182").Value = "Choclero"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 3000
$ws.Range("K182").Value = 420
$ws.Range("L182").Value = 450
$ws.Range("M182").Value = 435
$ws.Range("N182").Value = "$/unidad"
$ws.Range("O182").Value = "Región Metropolitana"
$ws.Range("P182").Value = 435
$ws.Range("Q182").Value = 1
$ws.Range("R182").Value = "Hortaliza"
